# TC19_Canine_Filter_Breed-CockerSpan.xlsx — "Fixed ICDC breed all testcases"
#
# The StatQuery column (C) on the "startup" sheet previously held the old
# "all breeds / all studies / all sexes" Cypher query in C2:C4. That query
# is replaced everywhere by a new Programs/Studies/Cases/Samples/Files
# rollup query. Everything else (TabName, query, dbExcel, WebExcel columns)
# stays the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$newStatQuery = 'MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)' + $nl +
    'OPTIONAL MATCH (samp:sample)-->(c)' + $nl +
    'OPTIONAL MATCH (diag:diagnosis)-->(c)' + $nl +
    'OPTIONAL MATCH (f:file)-[*]->(c)' + $nl +
    'OPTIONAL MATCH (sf:file)-->(s)' + $nl +
    'WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p' + $nl +
    "WHERE demo.breed IN ['Cocker Spaniel']" + $nl +
    'RETURN  ' + $nl +
    '    count(distinct p) AS Programs,' + $nl +
    '    count(distinct s) AS Studies,' + $nl +
    '    count(distinct c) AS Cases,' + $nl +
    '    count(distinct samp) AS Samples,' + $nl +
    '    count(distinct f) AS `Case Files`,' + $nl +
    '    count(distinct sf) AS `Study Files`'

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# View state: scrolled down so row 3 is at top, zoomed to 85%, B4 selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$win.Zoom = 85
$ws.Range("B4").Select()
